# This workbook ("Pais" sheet) is a daily Covid-19 country stats snapshot.
# The update (a) refreshes the "last updated" footer timestamp and (b) refreshes
# same-day totals for a handful of countries -- a few of which also needed their
# row reordered/relabelled so the country name lines up with its row of figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: footer timestamp "...16:22" -> "...16:52"
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Abril de 2020 a las 16:52"

# Row 4: updated figures
$ws.Cells.Item(4, 2).Value = 819443
$ws.Cells.Item(4, 3).Value = 699
$ws.Cells.Item(4, 4).Value = 83017
$ws.Cells.Item(4, 5).Value = 691055
$ws.Cells.Item(4, 7).Value = 53
$ws.Cells.Item(4, 8).Value = 45371

# Row 8: updated figures
$ws.Cells.Item(8, 5).Value = 44398
$ws.Cells.Item(8, 7).Value = 41
$ws.Cells.Item(8, 8).Value = 5127

# Row 64: updated figures
$ws.Cells.Item(64, 2).Value = 2070
$ws.Cells.Item(64, 3).Value = 75
$ws.Cells.Item(64, 4).Value = 515
$ws.Cells.Item(64, 5).Value = 1536

# Row 85: updated figures
$ws.Cells.Item(85, 2).Value = 1024
$ws.Cells.Item(85, 3).Value = 49
$ws.Cells.Item(85, 5).Value = 801
$ws.Cells.Item(85, 6).Value = 37
$ws.Cells.Item(85, 7).Value = 4
$ws.Cells.Item(85, 8).Value = 49

# Row 141: updated figures
$ws.Cells.Item(141, 4).Value = 37
$ws.Cells.Item(141, 5).Value = 70

# Row 142: updated figures
$ws.Cells.Item(142, 4).Value = 20
$ws.Cells.Item(142, 5).Value = 73

# Row 147: 'Maldivas' -> 'Togo'
$ws.Cells.Item(147, 1).Value = "Togo"
$ws.Cells.Item(147, 2).Value = 88
$ws.Cells.Item(147, 3).Value = 2
$ws.Cells.Item(147, 4).Value = 56
$ws.Cells.Item(147, 5).Value = 26
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 8).Value = 6

# Row 148: 'Togo' -> 'Maldivas'
$ws.Cells.Item(148, 1).Value = "Maldivas"
$ws.Cells.Item(148, 4).Value = 16
$ws.Cells.Item(148, 5).Value = 70
$ws.Cells.Item(148, 6).Value = 2
$ws.Cells.Item(148, 8).Value = 0

# Row 155: 'Islas Caimanes' -> 'Guyana'
$ws.Cells.Item(155, 1).Value = "Guyana"
$ws.Cells.Item(155, 2).Value = 67
$ws.Cells.Item(155, 3).Value = 1
$ws.Cells.Item(155, 4).Value = 9
$ws.Cells.Item(155, 5).Value = 51
$ws.Cells.Item(155, 6).Value = 5
$ws.Cells.Item(155, 8).Value = 7

# Row 156: 'Guyana' -> 'Islas Caimanes'
$ws.Cells.Item(156, 1).Value = "Islas Caimanes"
$ws.Cells.Item(156, 4).Value = 7
$ws.Cells.Item(156, 5).Value = 58
$ws.Cells.Item(156, 6).Value = 3
$ws.Cells.Item(156, 8).Value = 1

# Row 168: 'Puerto Rico' -> 'Mozambique'
$ws.Cells.Item(168, 1).Value = "Mozambique"
$ws.Cells.Item(168, 2).Value = 41
$ws.Cells.Item(168, 3).Value = 2
$ws.Cells.Item(168, 4).Value = 8
$ws.Cells.Item(168, 5).Value = 33
$ws.Cells.Item(168, 8).Value = 0

# Row 169: 'Eritrea' -> 'Puerto Rico'
$ws.Cells.Item(169, 1).Value = "Puerto Rico"
$ws.Cells.Item(169, 4).Value = 1
$ws.Cells.Item(169, 5).Value = 36
$ws.Cells.Item(169, 8).Value = 2

# Row 170: 'Mozambique' -> 'Eritrea'
$ws.Cells.Item(170, 1).Value = "Eritrea"
$ws.Cells.Item(170, 4).Value = 6
$ws.Cells.Item(170, 5).Value = 33

# Row 195: 'Islas Turcas y Caicos' -> 'Montserrat'
$ws.Cells.Item(195, 1).Value = "Montserrat"
$ws.Cells.Item(195, 4).Value = 2
$ws.Cells.Item(195, 5).Value = 9
$ws.Cells.Item(195, 6).Value = 1
$ws.Cells.Item(195, 8).Value = 0

# Row 196: 'Montserrat' -> 'Islas Malvinas'
$ws.Cells.Item(196, 1).Value = "Islas Malvinas"
$ws.Cells.Item(196, 4).Value = 3
$ws.Cells.Item(196, 5).Value = 8
$ws.Cells.Item(196, 6).Value = 0

# Row 197: 'Islas Malvinas' -> 'Burundi'
$ws.Cells.Item(197, 1).Value = "Burundi"
$ws.Cells.Item(197, 3).Value = 6
$ws.Cells.Item(197, 4).Value = 4
$ws.Cells.Item(197, 5).Value = 6
$ws.Cells.Item(197, 8).Value = 1

# Row 198: 'Seychelles' -> 'Islas Turcas y Caicos'
$ws.Cells.Item(198, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(198, 4).Value = 4
$ws.Cells.Item(198, 8).Value = 1

# Row 199: 'Burundi' -> 'Seychelles'
$ws.Cells.Item(199, 1).Value = "Seychelles"
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 5
$ws.Cells.Item(199, 8).Value = 0

# Row 215: 'Yemen' -> 'San Pedro y Miquelon'
$ws.Cells.Item(215, 1).Value = "San Pedro y Miquelon"

# Row 216: 'San Pedro y Miquelon' -> 'Yemen'
$ws.Cells.Item(216, 1).Value = "Yemen"
